$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.648594856262207
$ws.Range("B1").Value = 0.7740026116371155
$ws.Range("C1").Value = 0.8881377577781677
$ws.Range("D1").Value = 5.082324504852295
$ws.Range("E1").Value = 1.616735816001892
